$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Body Length [m]
$ws.Range("A2").Value = "Body Length [m]"
$ws.Range("B2").Value = 0.75

# Row 3: Body CG [m from front] (formula = B2/2)
$ws.Range("A3").Value = "Body CG [m from front]"
$ws.Range("B3").Formula = "=B2/2"

# Row 4: Thigh Distance [cm]
$ws.Range("A4").Value = "Thigh Distance [cm]"
$ws.Range("B4").Value = 35

# Row 5: Thigh Inertia [m4]
$ws.Range("A5").Value = "Thigh Inertia [m4]"
$ws.Range("B5").Value = 0.00067929

# Row 6: Thigh CG (from hip) [m]
$ws.Range("A6").Value = "Thigh CG (from hip) [m]"
$ws.Range("B6").Value = 0.2

# Row 7: Thigh Mass [kg]
$ws.Range("A7").Value = "Thigh Mass [kg]"
$ws.Range("B7").Value = 0.43

# Row 8: Shank Distance [cm]
$ws.Range("A8").Value = "Shank Distance [cm]"
$ws.Range("B8").Value = 35

# Row 9: Shank Inertia [m4]
$ws.Range("A9").Value = "Shank Inertia [m4]"
$ws.Range("B9").Value = 0.00051042

# Row 10: Shank CG (from knee) [m]
$ws.Range("A10").Value = "Shank CG (from knee) [m]"
$ws.Range("B10").Value = 0.19

# Row 11: Shank Mass [kg]
$ws.Range("A11").Value = "Shank Mass [kg]"
$ws.Range("B11").Value = 0.37
